$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.273867249488831
$ws.Range("B1").Value = 1.84672749042511
$ws.Range("C1").Value = 4.265507221221924
$ws.Range("D1").Value = 3.013212919235229
$ws.Range("E1").Value = 1.15052318572998
